# Generate Report for Handback
# Updates handback status timestamps/priority across the Overview, zh-cn and
# de-de sheets. Several cells across sheets reuse the same text value, so all
# occurrences of each changed value are updated together.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" for
# a5540b81-...: 2016-09-05 02:18:06 -> 2016-09-05 02:18:58
$wsOverview.Range("G4").Value = "2016-09-05 02:18:58"
$wsOverview.Range("G5").Value = "2016-09-05 02:18:58"
$wsDeDe.Range("H4").Value = "2016-09-05 02:18:58"
$wsDeDe.Range("H5").Value = "2016-09-05 02:18:58"

# Priority for a5540b81-...: ht -> mt
$wsZhCn.Range("E4").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"
$wsDeDe.Range("E4").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"

# zh-cn Correspond Handoff Datetime: 2016-09-05 02:17:58 -> 2016-09-05 02:18:53
$wsZhCn.Range("H4").Value = "2016-09-05 02:18:53"
$wsZhCn.Range("H5").Value = "2016-09-05 02:18:53"

# zh-cn Correspond Handback DateTime: 2016-09-05 02:18:28 -> 2016-09-05 02:19:16
$wsZhCn.Range("K4").Value = "2016-09-05 02:19:16"
$wsZhCn.Range("K5").Value = "2016-09-05 02:19:16"

# de-de Correspond Handback DateTime: 2016-09-05 02:18:35 -> 2016-09-05 02:19:23
$wsDeDe.Range("K4").Value = "2016-09-05 02:19:23"
$wsDeDe.Range("K5").Value = "2016-09-05 02:19:23"
